$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 92, pushing existing rows 92-93 down to 93-94.
# (Formatting, e.g. the date number format on column D, is inherited from the row above.)
$ws.Rows.Item(92).Insert()

# Populate the newly inserted row 92 with the new weekly price entry.
$ws.Cells.Item(92, 1).Value = 8
$ws.Cells.Item(92, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(92, 3).Value = "Coquimbo"
$ws.Cells.Item(92, 4).Value = 44509
$ws.Cells.Item(92, 5).Value = 4
$ws.Cells.Item(92, 6).Value = 100112044
$ws.Cells.Item(92, 7).Value = "Perejil"
$ws.Cells.Item(92, 8).Value = "Sin especificar"
$ws.Cells.Item(92, 9).Value = "Primera"
$ws.Cells.Item(92, 10).Value = 2940
$ws.Cells.Item(92, 11).Value = 1300
$ws.Cells.Item(92, 12).Value = 1500
$ws.Cells.Item(92, 13).Value = 1400
$ws.Cells.Item(92, 14).Value = "`$/atado 1 a 1,5 kilos"
$ws.Cells.Item(92, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(92, 16).Value = 933
$ws.Cells.Item(92, 17).Value = 1.5
$ws.Cells.Item(92, 18).Value = "Hortaliza"
